# Apply the cryptos-list refresh captured in the commit diff.
# Updates Price (D) and Volume(1h) (E) for every row, and for the
# handful of rows whose rank order changed, also updates Coin (B)
# and Link (C) so each row keeps the correct name/url/price/volume
# grouping while the row (and its index in column A) stays put.
#
# All target cells hold plain text (inline/shared strings) with no
# explicit cell style in the source workbook. Several of the new
# Price values (e.g. "1.000", "5.230", "1.510") are valid numeric
# literals, so a bare .Value assignment would let Excel parse them
# as numbers and silently drop meaningful trailing zeros. To avoid
# that, every write here is done as: force text format ("@"),
# assign the literal string, then reset the style back to Normal
# so the cell ends up with the original (default, unstyled) look
# while still holding the exact text from the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "30.389.94"
Set-TextValue "E2" "  +0.76%  "
Set-TextValue "D3" "1.880.02"
Set-TextValue "E3" "  +0.97%  "
Set-TextValue "D4" "1.000"
Set-TextValue "E4" "  -0.07%  "
Set-TextValue "D5" "245.11"
Set-TextValue "E5" "  +4.74%  "
Set-TextValue "D6" "0.9996"
Set-TextValue "E6" "  -0.09%  "
Set-TextValue "D7" "0.4766"
Set-TextValue "E7" "  +1.77%  "
Set-TextValue "D8" "0.2875"
Set-TextValue "E8" "  +1.10%  "
Set-TextValue "E9" "  -0.36%  "
Set-TextValue "D10" "21.27"
Set-TextValue "E10" "  -0.38%  "
Set-TextValue "D11" "0.07755"
Set-TextValue "E11" "  +0.08%  "
Set-TextValue "D12" "1.893.93"
Set-TextValue "E12" "  +2.33%  "
Set-TextValue "D13" "96.61"
Set-TextValue "E13" "  +1.20%  "
Set-TextValue "D14" "0.7344"
Set-TextValue "E14" "  +6.85%  "
Set-TextValue "D15" "5.124"
Set-TextValue "E15" "  +0.94%  "
Set-TextValue "D16" "273.98"
Set-TextValue "E16" "  +3.42%  "
Set-TextValue "D17" "30.373.52"
Set-TextValue "E17" "  +0.73%  "
Set-TextValue "D18" "13.39"
Set-TextValue "E18" "  -1.58%  "
Set-TextValue "D19" "0.000007539"
Set-TextValue "E19" "  -2.47%  "
Set-TextValue "D20" "0.9999"
Set-TextValue "E20" "  -0.07%  "
Set-TextValue "D21" "2.131.28"
Set-TextValue "E21" "  -0.06%  "
Set-TextValue "D22" "0.9998"
Set-TextValue "E22" "  -0.05%  "
Set-TextValue "D23" "5.230"
Set-TextValue "E23" "  +0.20%  "
Set-TextValue "D24" "6.167"
Set-TextValue "E24" "  +0.29%  "
Set-TextValue "D25" "9.241"
Set-TextValue "E25" "  -2.09%  "
Set-TextValue "D26" "163.52"
Set-TextValue "E26" "  -1.47%  "
Set-TextValue "E27" "  +1.51%  "
Set-TextValue "D28" "1.958"
Set-TextValue "E28" "  +1.85%  "
Set-TextValue "B29" "Stellar"
Set-TextValue "C29" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D29" "0.09982"
Set-TextValue "E29" "  +0.84%  "
Set-TextValue "B30" "Toncoin"
Set-TextValue "C30" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D30" "1.368"
Set-TextValue "E30" "  +0.04%  "
Set-TextValue "D31" "1.510"
Set-TextValue "E31" "  +3.56%  "
Set-TextValue "D32" "4.306"
Set-TextValue "E32" "  -0.50%  "
Set-TextValue "D33" "4.074"
Set-TextValue "E33" "  +1.09%  "
Set-TextValue "D34" "0.04745"
Set-TextValue "E34" "  +0.64%  "
Set-TextValue "E35" "  -0.22%  "
Set-TextValue "D36" "0.6956"
Set-TextValue "E36" "  +0.00%  "
Set-TextValue "D37" "2.716"
Set-TextValue "E37" "  -0.14%  "
Set-TextValue "D38" "0.01859"
Set-TextValue "E38" "  +0.20%  "
Set-TextValue "D39" "2.750"
Set-TextValue "E39" "  -0.40%  "
Set-TextValue "D40" "6.270"
Set-TextValue "E40" "  -0.39%  "
Set-TextValue "D41" "0.8423"
Set-TextValue "E41" "  +1.29%  "
Set-TextValue "D42" "69.31"
Set-TextValue "E42" "  -3.34%  "
Set-TextValue "B43" "RenderToken"
Set-TextValue "C43" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D43" "1.905"
Set-TextValue "E43" "  -0.75%  "
Set-TextValue "B44" "TheSandbox"
Set-TextValue "C44" "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D44" "0.4161"
Set-TextValue "E44" "  +1.02%  "
Set-TextValue "B45" "PaxDollar"
Set-TextValue "C45" "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextValue "D45" "0.9991"
Set-TextValue "E45" "  -0.11%  "
Set-TextValue "D46" "101.85"
Set-TextValue "E46" "  -0.73%  "
Set-TextValue "B47" "Aptos"
Set-TextValue "C47" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue "D47" "7.079"
Set-TextValue "E47" "  -0.12%  "
Set-TextValue "B48" "EnergySwap"
Set-TextValue "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "9.198"
Set-TextValue "E48" "  +0.72%  "
Set-TextValue "E49" "  +1.63%  "
Set-TextValue "D50" "911.72"
Set-TextValue "E50" "  -5.69%  "
Set-TextValue "D51" "0.05593"
Set-TextValue "E51" "  -0.66%  "
